$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update room names to room numbers per the new room-name-to-number mapping
$ws.Range("B7").Value = "Private Lesson with Ivy CHUANG `n(Room G14)"
$ws.Range("E7").Value = "Private Lesson with Ivy CHUANG `n(Room G14)"
$ws.Range("F7").Value = "Flute MasterClass`n(Room G19)"

$ws.Range("B11").Value = "Rehearsal with pianist`n(Room G22)"
$ws.Range("D11").Value = "Private Lesson with Stephane RETY `n(Room G19)"

$ws.Range("C19").Value = "Private Lesson with Stephane RETY & pianist `n(Room G19)"
$ws.Range("F19").Value = "Flute MasterClass`n(Room G19)"

$ws.Range("D23").Value = "Ensemble `n(Room 242)"

# Restore automatic row heights (Excel recalculates wrap height on edit;
# AutoFit keeps rows at their natural/default height instead of pinning
# an explicit custom height).
$ws.Rows(7).AutoFit()
$ws.Rows(11).AutoFit()
$ws.Rows(19).AutoFit()
$ws.Rows(23).AutoFit()
